# Atualização de documentos. Correção de métricas no dashboard estático
#
# This script reproduces, via the PowerPoint COM object model, the layout
# fix applied to the single diagram slide: two icon groups ("Group 11" and
# "Group 10") were resized taller because one picture inside each group was
# nudged downward, and a couple of nearby shapes (a curved connector and a
# text box) were repositioned slightly to keep everything lined up.
#
# Because PowerPoint's object model stores shape position/size in points
# while the underlying OOXML stores EMU (1 pt = 12700 EMU), and because the
# group's child-coordinate-space (chOff/chExt) is only recomputed by the
# host when a group is actually re-grouped, the safe way to reproduce the
# exact target geometry is to: ungroup -> reposition the children -> regroup.
# Re-grouping also naturally makes chOff/chExt collapse onto off/ext, which
# is exactly what the target file shows for "Group 11".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

# Burns one shape-id from the slide's id allocator (PowerPoint never reuses
# shape ids, even after the shape is deleted) so that a subsequent Group()
# call lands on a specific, predictable id.
function Burn-ShapeId() {
    $tmp = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
    $tmp.Delete()
}

# ---------------------------------------------------------------------
# 1) "Group 11" (id 12) -> becomes "Group 2" (id 3)
#    The laptop picture (id 7) moves down/left a bit inside the group,
#    the call-center picture (id 9) becomes the new child-space origin,
#    and the group's own box grows from 914400 to 1007164 EMU tall.
# ---------------------------------------------------------------------
Burn-ShapeId    # consumes id 2, so the regroup below lands on id 3

$group1 = Get-ShapeById $s.Shapes 12
$items1 = $group1.Ungroup()

$laptop = $null
$callCenter1 = $null
for ($i = 1; $i -le $items1.Count; $i++) {
    $sh = $items1.Item($i)
    if ($sh.Id -eq 7) { $laptop = $sh }
    if ($sh.Id -eq 9) { $callCenter1 = $sh }
}

$laptop.Left = 124.157956
$laptop.Top = 142.494957
$callCenter1.Left = 62.592757
$callCenter1.Top = 135.190712

$range1 = $s.Shapes.Range(@($laptop.Name, $callCenter1.Name))
$newGroup1 = $range1.Group()
$newGroup1.Name = "Group 2"

# ---------------------------------------------------------------------
# 2) "Group 10" (id 11) -> becomes "Group 5" (id 6)
#    The computer picture (id 5) moves straight down inside the group,
#    the group's own box grows from 914400 to 1007164 EMU tall.
# ---------------------------------------------------------------------
Burn-ShapeId    # consumes id 4, so the regroup below lands on id 6

$group2 = Get-ShapeById $s.Shapes 11
$items2 = $group2.Ungroup()

$computer = $null
$callCenter2 = $null
for ($i = 1; $i -le $items2.Count; $i++) {
    $sh = $items2.Item($i)
    if ($sh.Id -eq 5) { $computer = $sh }
    if ($sh.Id -eq 10) { $callCenter2 = $sh }
}

$computer.Left = 161.201424
$computer.Top = 70.494961

$range2 = $s.Shapes.Range(@($computer.Name, $callCenter2.Name))
$newGroup2 = $range2.Group()
$newGroup2.Name = "Group 5"

# ---------------------------------------------------------------------
# 3) Curved connector between the laptop group and the database icon
#    shifts slightly and gets a bit narrower.
# ---------------------------------------------------------------------
$connector = Get-ShapeById $s.Shapes 28
$connector.Left = 167.743310
$connector.Top = 206.909600
$connector.Width = 92.062596
$connector.Height = 107.233227

# ---------------------------------------------------------------------
# 4) The "Armazenamos ..." text box shifts slightly left/down to stay
#    aligned under its icon.
# ---------------------------------------------------------------------
$textBox = Get-ShapeById $s.Shapes 41
$textBox.Left = 218.741189
$textBox.Top = 354.332535
